$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: binomial distribution -> geometric distribution
$ws.Range("D5").Value = "기하 분포"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/04/28/geometric_distribution.html"

# Row 9: Neural Net article -> foreign degree job market article
$ws.Range("D9").Value = "해외대학이 대기업 취직에 유리한 이유 (3)"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/foreign-degree-job-market-merits-3/#utm_source=rss&utm_medium=rss&utm_campaign=foreign-degree-job-market-merits-3"

# Row 51: datetime formatting article -> pandas excel article
$ws.Range("D51").Value = "[python] pandas로 엑셀 파일 읽고 수정한 후 쓰기"
$ws.Range("E51").Value = "https://bskyvision.com/1183"
